# Auto-generated PowerShell COM-interop script
# Applies the 'artificer tier 2 and 3' edit to Skills_Table (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at I, shifting the old I (boolean 'Spell' flag) to J.
#    This preserves all existing I-column values/styles, moving them to J,
#    and leaves a blank (but styled, for row 1) column I behind.
$ws.Columns("I").Insert()

# 2) Header row: new I1 = 'Special' (J1 retains 'Spell' from the shift above).
$ws.Range("I1").Value = "Special"

# 3) Targeted content updates on existing rows 189-195.
$ws.Range("A189").Value = 'Appraise [Constructs]'
$ws.Range("A191").Value = 'Appraise [Wondrous Items]'

$ws.Range("B192").Value = '[Wondrous Items] You may attach a mechanism to an object that alters its properties in one way of your choice with no gem cost: the object (1) sheds light, as though through a light spell, (2) shows a static visual effect (or small picture), (3) makes a continuous non-damaging audio effect or a nonverbal sound, or (4) plays a recorded message when touched, up to 6 seconds long.'
$ws.Range("B193").Value = '[Runes] You may disguise the meaning of your Artificer Schematics, so other fellow Artificers cannot read it to learn its secrets.  Upon creation of a Schematic, you may opt to encrypt it.  The Professional expertise required to bypass your encryption for this Schematic is your Artificer Professional Expertise Tier when the Schematic was created.'
$ws.Range("B195").Value = '[Weaponry, Alchemy] Through careful art, science, and ritual, you are able to silver-plate to your weapon’s blade like Sir Artorias the Moonslayer.  Until the end of the day, your blade is silvered, allowing your attacks with this weapon to be harrowing against lycanthropes, wraiths, vampires, and other cursed entities.'

# Row 193: I193 gains new augment text (was blank).
$ws.Range("I193").Value = 'Characters with Appraise [Runes] may, once per event, use a Fusion Point to increase the tier to bypass encryption (see below) by 1.  This increase may be countered by knowledge of how the runes work (the creator innately ignores it) or extracting the information from the creator.'

# Row 195: G195 gains new prerequisite text (was blank).
$ws.Range("G195").Value = 'Appraise [Alchemy]'

# 4) Append new rows 196-207 (Artificer tier 2 and 3 skills).
$newRows = New-Object 'object[,]' 12,10
$newRows[0,0] = 'Bypass Encryption 1'
$newRows[0,1] = '[Runes]. You are able to bypass other Artificers’ Schematic Encryption.  Spend 15 minutes per Tier of the encryption to break the code, being able to use this as a reference material for use to create your own Schematic.'
$newRows[0,2] = 'Artificer'
$newRows[0,3] = 2
$newRows[0,4] = 'You must be In-Play for the time to count.  If used as a reference material for replicating a Schematic for yourself, you mark up the formerly encrypted Schematic so much that it ceases to be useful and is effectively consumed.  Reference materials of decrypted Schematics cannot be enhanced by your own Schematic encryption at this stage.  It was difficult enough to decrypt.  You may only bypass an encryption once per event.'
$newRows[0,5] = ""
$newRows[0,6] = ""
$newRows[0,7] = ""
$newRows[0,8] = 'Characters with Appraise [Runes] may, once per event, use a Fusion Point to increase their ability to bypass encryption by 1.'
$newRows[0,9] = $false
$newRows[1,0] = 'Artificer’s Wayfinder'
$newRows[1,1] = '[Wondrous Items, Clothing]. You create a curious compass or pathfinding accessory that makes it easy to traverse the Wilderness for this event.  One adventuring party of your choice who wields the wayfinder can treat Forest and Blasted Lands hexes as 1 Movement Point, rather than 2.'
$newRows[1,2] = 'Artificer'
$newRows[1,3] = 2
$newRows[1,4] = 'This item counts toward the limit of equipped Artificer Objects a character may attune to at once. Multiple compasses in an adventuring party do not stack.  This Wayfinder requires an Artificer Object, but no schematic (unless you want to).'
$newRows[1,5] = 'A compass, pocket-watch, or ornate clothing accessory.  An Object Schematic to match.'
$newRows[1,6] = 'Tinkerer’s Quirk'
$newRows[1,7] = 'You may consume a Small gem and a Fusion Point to increase your adventuring party’s maximum movement points by 2.'
$newRows[1,8] = 'If you have Appraise [Wondrous Items] or Appraise [Clothing] and augment the wayfinder, reduce all Movement point costs by 1, including Mountains (minimum'
$newRows[1,9] = $false
$newRows[2,0] = 'Alchemical Concoction 1'
$newRows[2,1] = '[Alchemy, Armor, Clothing]. You tinker with your headgear, experimenting with protective oils, salves, and bulwarks that bolster your safety– at least at first blush.  Once per event, you may spend two minutes applying a concoction to headgear.  This headgear confers a one-time additional floating hit. This bonus disappears once struck; it cannot be regenerated until next event.'
$newRows[2,2] = 'Artificer'
$newRows[2,3] = 2
$newRows[2,4] = 'This weapon counts toward the limit of equipped Artificer Objects a character may attune to at once.'
$newRows[2,5] = 'Relevant headgear, such as a helm or hat.'
$newRows[2,6] = 'Appraise [Alchemy]'
$newRows[2,7] = ""
$newRows[2,8] = 'You may take this skill up to twice, once for each kind of headgear (Armor or Clothing).'
$newRows[2,9] = $false
$newRows[3,0] = 'Weapon Warp'
$newRows[3,1] = '[Weaponry].  You may alter a weapon’s properties with trace energy so that it feels easy to wield.  For one encounter, a wielder of this weapon can ignore their base class weapon restrictions.'
$newRows[3,2] = 'Artificer'
$newRows[3,3] = 2
$newRows[3,4] = 'This weapon counts toward the limit of equipped Artificer Objects a character may attune to at once.  This does not make your current one-handed weapon deal damage as though it is two-handed.  You simply can utilize a two-handed one, for example, and it doesn’t penalize you (i.e., you physically use a greatsword, not pretend like a short-sword is a greatsword).  Skills with specific limitations (e.g., Rogue’s Backstab) still require specific weapons to be utilized.'
$newRows[3,5] = 'A weapon.  No schematic is necessary (unless you want to add to your list of schematics).'
$newRows[3,6] = ""
$newRows[3,7] = 'You may use a Small gem and a Fusion Point to extend this duration to four hours.'
$newRows[3,8] = ""
$newRows[3,9] = $false
$newRows[4,0] = 'Alchemist’s Fire'
$newRows[4,1] = '[Alchemy, Weaponry]. Through curious alchemical tinkering, you’ve bottled congealed fire.  Expend one Fusion Point.  You gain one vial of alchemist fire.  You can hurl it or use it to make your weapon wreath itself in flame.   Hurling alchemist fire deals 1 Fire Damage upon impact, which bypasses shields (but not armor).'
$newRows[4,2] = 'Artificer'
$newRows[4,3] = 2
$newRows[4,4] = 'Only a character with the Alchemist’s Fire Skill may safely use Alchemist’s Fire (projectile or weapon) without instantly going up in flames themselves.'
$newRows[4,5] = 'Thrown alchemist fire should mirror a LARP-safe fireball or nerf-football.  Applying alchemist fire to a weapon should be represented by a red ribbon. Attacking should be telegraphed by stating “Alchemist Fire!” then stating the damage (“1”), “Fire damage,” and “through Shields” if thrown.'
$newRows[4,6] = 'Alchemical Concoction 1.'
$newRows[4,7] = ""
$newRows[4,8] = 'Having Appraise [Weapon] enables an additional option'
$newRows[4,9] = $false
$newRows[5,0] = 'Dedicated Artisan'
$newRows[5,1] = 'Focusing on their Profession, the Artificer begins to derive concentrated benefit from their creations.  While all other characters may attune to only one Artificer Object, the Artisan Artificer may don one additional Artificer Object, now benefiting from up to two effects.'
$newRows[5,2] = 'Artificer'
$newRows[5,3] = 3
$newRows[5,4] = ""
$newRows[5,5] = ""
$newRows[5,6] = ""
$newRows[5,7] = ""
$newRows[5,8] = ""
$newRows[5,9] = $false
$newRows[6,0] = 'Keen Eye'
$newRows[6,1] = '[Gems].  You have a practiced eye for leveraging Gems and unlocking outsized potential.  Once per event, you may increase a Gem’s value by one level (see The Object’s Power Source below) when consuming it in the following ways: to power an Artificer skill, serve as a Power Source of an Object, or trade it to Organizers for Faction Gold Dragons.'
$newRows[6,2] = 'Artificer'
$newRows[6,3] = 3
$newRows[6,4] = 'The appraised Gem may not exceed your Artificer Tier in value.  For example, as a Tier 3 Artificer, you may increase one Gem from Medium (Tier 2) to Large (Tier 3), but not from Large (Tier 3) to Exquisite (Tier 4).'
$newRows[6,5] = ""
$newRows[6,6] = 'Appraise [Gems]'
$newRows[6,7] = ""
$newRows[6,8] = ""
$newRows[6,9] = $false
$newRows[7,0] = 'Craft Homunculus'
$newRows[7,1] = '[Constructs, Wondrous Items]. You forge a helper-automaton to aid you in your Artificer tasks through a ritual with an audience.  At the end of the ritual, expend a Fusion Point.  When a homunculus is on your person and engaged in a scene, you may reduce your character’s time to attune to an Artificer Object by 1 minute (minimum 1 minute).  Once ‘built,’ the homunculus persists indefinitely and does not apply to the limit of Artificer items created in future events.'
$newRows[7,2] = 'Artificer'
$newRows[7,3] = 3
$newRows[7,4] = 'A plush doll, stuffed animal, figurine, 3d printed creature, or something else along those lines—tailored to look like a clockwork or created creature.  This should be plainly visible in any scene in which it is used and LARP-safe. Limitations'
$newRows[7,5] = 'A plush doll, stuffed animal, figurine, 3d printed creature, or something else along those lines—tailored to look like a clockwork or created creature.  This should be plainly visible in any scene in which it is used and LARP-safe. Limitations'
$newRows[7,6] = 'Appraise [Constructs] or Appraise [Wondrous Items]'
$newRows[7,7] = 'You may expend a Tiny gem and a Fusion Point to increase the reduction to 2 minutes (minimum 1 minute) for 4 hours.  You may expend a Small Gem and a Fusion to reduce it by 2 minutes (minimum 1 minute) AND grant your Homunculus the use of [Cooperative Action] for purposes of encrypting or decrypting Schematics for 4 hours.'
$newRows[7,8] = 'You can increase a homunculus’ defense by casting Arcane Armor on it or by making the homunculus an Artificer Object with an Arcane Armor Schematic.  Knowing both Appraise [Constructs] AND Appraise [Wondrous Items] additionally prevents the doubled time penalty for [Cooperative Action] with your homunculus.'
$newRows[7,9] = $false
$newRows[8,0] = 'Tailor’s Threads'
$newRows[8,1] = '[Clothing]. A brilliant clothier can make one feel lighter on their feet.  For one encounter, you may alter a clothing set’s properties.  The wearer’s maximum AP limitation increases by 1 (to a maximum of 4) while wearing this armor.'
$newRows[8,2] = 'Artificer'
$newRows[8,3] = 3
$newRows[8,4] = 'This is limited to outfits consisting primarily of cloth, hide, and leather.  This armor counts toward the limit of equipped Artificer Objects a character may attune to at once.  This does not give you any bonus to AP.  It simply allows you to wear heavier armors without penalty (for purposes of spellcasting, sneaking, etc.) for a temporary period.'
$newRows[8,5] = 'A fancy outfit'
$newRows[8,6] = 'Appraise [Clothing]'
$newRows[8,7] = 'You may consume a Large gem and a Fusion Point to extend this effect to 4 hours.'
$newRows[8,8] = ""
$newRows[8,9] = $false
$newRows[9,0] = 'Armor Infusion'
$newRows[9,1] = '[Armor]. In the fires of a forge, you make the most of a piece of armor. For one encounter, you may alter an armor’s properties to make it denser and more protective.  The armor’s AP value to one hit zone increases by 1 (to a maximum of 4) while wearing this armor.'
$newRows[9,2] = 'Artificer'
$newRows[9,3] = 3
$newRows[9,4] = 'This is limited to AP1+ armors; it becomes heavier, which may limit certain Path abilities. This armor counts toward the limit of equipped Artificer Objects a character may attune to at once.'
$newRows[9,5] = 'Armor, gauntlets, greaves, etc.'
$newRows[9,6] = 'Alchemical Concoction 1'
$newRows[9,7] = 'You may consume a Large gem and a Fusion Point to extend this effect to 4 hours.'
$newRows[9,8] = ""
$newRows[9,9] = $false
$newRows[10,0] = 'Riastrad’s Magic Traps'
$newRows[10,1] = '[Wondrous Items, Alchemy]. Your adept fabrication of spell-like effects allows you to create unique approaches to traps, not unlike Riastrad the Fireforged.  You may apply effects to existing traps that mimic the consequences of spells for which you have Object Schematics (see Creating an Artificer Object below).  If triggered, the trap activates this spell-like ability instead of a poison.'
$newRows[10,2] = 'Artificer'
$newRows[10,3] = 3
$newRows[10,4] = 'A character who has Trap Use 1 or Area Trap—it need not be you, but you must create and set it together in a 3 minute RP scene.  A Riastrad’s Magic Trap is considered an Artificer Object and requires Fusion Points to bind the spell-like ability to the trap (see Eligible Skills/Spells to Artifice, below).  It is not considered ‘equipped’ like an Artificer Object is.  However, it deducts from the number of mundane Traps a character is able to create per event.  The effect can be dispelled with Break Arcana or similar alchemical abilities that would remove Artificer Effects.   If the trap is sprung or destroyed (per the T3 Rogue Skill, but NOT T4), the effect dissipates.'
$newRows[10,5] = 'A typical Trap that allows Rogue workflow from Detect Traps to Disarm Traps, as well as has the reference material’s effects if triggered.'
$newRows[10,6] = 'Appraise [Wondrous Items] or Appraise [Alchemy]'
$newRows[10,7] = ""
$newRows[10,8] = ""
$newRows[10,9] = $false
$newRows[11,0] = 'Fabricate Harvest Guard'
$newRows[11,1] = '[Construct]. From a solution of straw, vegetation, and fear, you create an inert simulacrum of the Harvest Guard until the end of the encounter.  You may animate this construct with a Large Gem and 2 Fusion Points and bring it with you to a Land Search as an additional guardian bound to defend the lands you state.  The Harvest Guard construct is immune to Fear effects, poison, and piercing damage. In combat, it takes three hits to collapse (but takes double Fire damage).'
$newRows[11,2] = 'Artificer'
$newRows[11,3] = 3
$newRows[11,4] = 'This circlet counts toward the limit of equipped Artificer Objects a character may attune to at once.  You may only have one construct bound to you at a time.  Note that the Harvest Guard is given a directive to defend the lands, not to defend a specific person—it will not willingly leave this hex unless you create another.  Depending on the commands given, losing control of the circlet may make the construct hostile toward the Artificer.'
$newRows[11,5] = 'A volunteer wearing a Harvest Guard costume (Artificer provides).  Harvest Guard constructs are eligible to use curved blades, like scythes, of any length.  The Artificer wears a lootable circlet that binds the construct to their command. The Artificer must create an Object Schematic reflective of the Harvest Guard.  This requires having met and studied the Harvest Guard (effectively Appraising them) for multiple events or relied on an expert’s Harvest Guard knowledge to create one.  A 10 minute RP scene.'
$newRows[11,6] = 'Craft Homunculus'
$newRows[11,7] = ""
$newRows[11,8] = ""
$newRows[11,9] = $false

$ws.Range("A196:J207").Value = $newRows

# Done.
